$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'50.762.38"
$ws.Range("E2").Value = "'  -2.75%  "
$ws.Range("D3").Value = "'2.728.67"
$ws.Range("E3").Value = "'  -3.17%  "
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'349.79"
$ws.Range("E5").Value = "'  -2.93%  "
$ws.Range("D6").Value = "'106.24"
$ws.Range("E6").Value = "'  -4.34%  "
$ws.Range("D7").Value = "'0.542"
$ws.Range("E7").Value = "'  -4.00%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "'  -0.08%  "
$ws.Range("D9").Value = "'0.571"
$ws.Range("E9").Value = "'  -4.49%  "
$ws.Range("D10").Value = "'38.63"
$ws.Range("E10").Value = "'  -4.63%  "
$ws.Range("E11").Value = "'  +2.52%  "
$ws.Range("D12").Value = "'0.0824"
$ws.Range("E12").Value = "'  -4.09%  "
$ws.Range("D13").Value = "'19.39"
$ws.Range("E13").Value = "'  -1.91%  "
$ws.Range("D14").Value = "'7.38"
$ws.Range("D15").Value = "'3.168.04"
$ws.Range("E15").Value = "'  -2.89%  "
$ws.Range("D16").Value = "'2.782.83"
$ws.Range("E16").Value = "'  -3.33%  "
$ws.Range("D17").Value = "'0.908"
$ws.Range("E17").Value = "'  -1.14%  "
$ws.Range("D18").Value = "'50.766.46"
$ws.Range("E18").Value = "'  -2.57%  "
$ws.Range("D19").Value = "'7.59"
$ws.Range("E19").Value = "'  +1.45%  "
$ws.Range("D20").Value = "'2.99"
$ws.Range("E20").Value = "'  -4.37%  "
$ws.Range("D21").Value = "'12.79"
$ws.Range("E21").Value = "'  -3.87%  "
$ws.Range("E22").Value = "'  -4.66%  "
$ws.Range("D23").Value = "'68.75"
$ws.Range("E23").Value = "'  -1.91%  "
$ws.Range("D24").Value = "'261.13"
$ws.Range("E24").Value = "'  -4.41%  "
$ws.Range("D25").Value = "'2.67"
$ws.Range("E25").Value = "'  -4.72%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "'  +0.17%  "
$ws.Range("D27").Value = "'25.60"
$ws.Range("E27").Value = "'  -4.36%  "
$ws.Range("D28").Value = "'0.158"
$ws.Range("E28").Value = "'  +11.48%  "
$ws.Range("E29").Value = "'  +0.01%  "
$ws.Range("E30").Value = "'  -3.10%  "
$ws.Range("D31").Value = "'51.29"
$ws.Range("E31").Value = "'  -1.60%  "
$ws.Range("D32").Value = "'34.16"
$ws.Range("E32").Value = "'  -0.77%  "
$ws.Range("D33").Value = "'5.94"
$ws.Range("E33").Value = "'  +2.07%  "
$ws.Range("D34").Value = "'0.0435"
$ws.Range("E34").Value = "'  -8.83%  "
$ws.Range("D35").Value = "'0.0819"
$ws.Range("E35").Value = "'  -3.21%  "
$ws.Range("D36").Value = "'5.10"
$ws.Range("E36").Value = "'  -7.11%  "
$ws.Range("E37").Value = "'  +0.00%  "
$ws.Range("D38").Value = "'18.20"
$ws.Range("E38").Value = "'  +0.53%  "
$ws.Range("D39").Value = "'3.09"
$ws.Range("E39").Value = "'  -3.67%  "
$ws.Range("E40").Value = "'  -5.23%  "
$ws.Range("E41").Value = "'  -3.69%  "
$ws.Range("E42").Value = "'  -4.08%  "
$ws.Range("D43").Value = "'2.19"
$ws.Range("E43").Value = "'  -3.33%  "
$ws.Range("D44").Value = "'119.17"
$ws.Range("E44").Value = "'  -5.14%  "
$ws.Range("D45").Value = "'21.55"
$ws.Range("E45").Value = "'  -4.08%  "
$ws.Range("D46").Value = "'2.060.07"
$ws.Range("E46").Value = "'  -0.38%  "
$ws.Range("E47").Value = "'  -1.01%  "
$ws.Range("E48").Value = "'  -3.96%  "
$ws.Range("E49").Value = "'  -7.43%  "
$ws.Range("D50").Value = "'0.901"
$ws.Range("E50").Value = "'  -5.04%  "
$ws.Range("B51").Value = "'TrustWalletToken"
$ws.Range("C51").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.27"
$ws.Range("E51").Value = "'  +3.05%  "
